# Apply 2022-10-27 data update to violent-crime-full-year.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('I2').Value = 6021
$ws.Range('I3').Value = 6268
$ws.Range('E4').Value = 1970
$ws.Range('I4').Value = 1446
$ws.Range('I5').Value = 584
$ws.Range('I6').Value = 7111
$ws.Range('E7').Value = 25974
$ws.Range('I7').Value = 21430

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('I4').Value = 37
$ws.Range('I7').Value = 249

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('I3').Value = 54
$ws.Range('I4').Value = 17
$ws.Range('I6').Value = 87
$ws.Range('I7').Value = 231

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range('I6').Value = 21
$ws.Range('I7').Value = 73

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('I6').Value = 198
$ws.Range('I7').Value = 682

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('I2').Value = 107
$ws.Range('I3').Value = 143
$ws.Range('I7').Value = 387

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('I6').Value = 255
$ws.Range('I7').Value = 827

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('I4').Value = 11
$ws.Range('I7').Value = 213

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('I4').Value = 86
$ws.Range('I7').Value = 672
$ws.Range('I8').Value = 1285
$ws.Range('I9').Value = 106
$ws.Range('I12').Value = 49
$ws.Range('I15').Value = 247
$ws.Range('I19').Value = 594
$ws.Range('I20').Value = 534
$ws.Range('I23').Value = 214
$ws.Range('I24').Value = 60
$ws.Range('I25').Value = 114
$ws.Range('I29').Value = 1326
$ws.Range('I30').Value = 73
$ws.Range('I31').Value = 213
$ws.Range('I33').Value = 969
$ws.Range('I36').Value = 289
$ws.Range('I37').Value = 682
$ws.Range('I42').Value = 742
$ws.Range('I44').Value = 159
$ws.Range('E46').Value = 57
$ws.Range('I48').Value = 287
$ws.Range('I50').Value = 106
$ws.Range('I51').Value = 246
$ws.Range('I54').Value = 438
$ws.Range('I55').Value = 237
$ws.Range('I60').Value = 117
$ws.Range('I63').Value = 78
$ws.Range('I67').Value = 827
$ws.Range('I70').Value = 34
$ws.Range('I73').Value = 197
$ws.Range('I76').Value = 308
$ws.Range('I78').Value = 292
$ws.Range('I79').Value = 607
$ws.Range('I85').Value = 975
$ws.Range('I86').Value = 131
$ws.Range('I89').Value = 249
$ws.Range('I90').Value = 266
$ws.Range('I94').Value = 223
$ws.Range('I96').Value = 231
$ws.Range('I98').Value = 147
$ws.Range('I99').Value = 387
$ws.Range('E101').Value = 25974
$ws.Range('I101').Value = 21430

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('I2').Value = 216
$ws.Range('I3').Value = 368
$ws.Range('I4').Value = 42
$ws.Range('I7').Value = 969

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('I3').Value = 97
$ws.Range('I6').Value = 209
$ws.Range('I7').Value = 438

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('I2').Value = 388
$ws.Range('I3').Value = 457
$ws.Range('I7').Value = 1326

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('I2').Value = 201
$ws.Range('I6').Value = 176
$ws.Range('I7').Value = 594

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('I2').Value = 50
$ws.Range('I7').Value = 159

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('I4').Value = 36
$ws.Range('I7').Value = 287

$ws = $wb.Worksheets.Item('River North')
$ws.Range('I2').Value = 59
$ws.Range('I7').Value = 308

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('I3').Value = 376
$ws.Range('I7').Value = 975

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('I2').Value = 186
$ws.Range('I5').Value = 25
$ws.Range('I7').Value = 742

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('I3').Value = 29
$ws.Range('I6').Value = 67

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('I2').Value = 68
$ws.Range('I7').Value = 292

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('I2').Value = 71
$ws.Range('I3').Value = 75
$ws.Range('I7').Value = 237

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range('I2').Value = 21
$ws.Range('I7').Value = 60

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range('E4').Value = 8
$ws.Range('E7').Value = 57

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('I6').Value = 62
$ws.Range('I7').Value = 214

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('I3').Value = 196
$ws.Range('I6').Value = 179
$ws.Range('I7').Value = 607

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('I4').Value = 33
$ws.Range('I6').Value = 187
$ws.Range('I7').Value = 534

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('I3').Value = 94
$ws.Range('I7').Value = 289

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('I3').Value = 165
$ws.Range('I4').Value = 38

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('I6').Value = 127
$ws.Range('I7').Value = 223

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('I6').Value = 30
$ws.Range('I7').Value = 114

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('I3').Value = 56
$ws.Range('I7').Value = 247

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('I6').Value = 96
$ws.Range('I7').Value = 147

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('I6').Value = 31
$ws.Range('I7').Value = 106

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range('I3').Value = 37
$ws.Range('I7').Value = 106

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('I2').Value = 64
$ws.Range('I7').Value = 197

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range('I3').Value = 9
$ws.Range('I7').Value = 34

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('I3').Value = 363
$ws.Range('I4').Value = 80
$ws.Range('I6').Value = 417
$ws.Range('I7').Value = 1285

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('I6').Value = 32
$ws.Range('I7').Value = 131

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('I2').Value = 89
$ws.Range('I6').Value = 90
$ws.Range('I7').Value = 266

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('I4').Value = 26
$ws.Range('I7').Value = 246

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('I2').Value = 41
$ws.Range('I7').Value = 117

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('I2').Value = 221
$ws.Range('I6').Value = 177
$ws.Range('I7').Value = 672

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range('I6').Value = 29
$ws.Range('I7').Value = 86

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range('I3').Value = 8
$ws.Range('I7').Value = 49
